$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.615.63"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.289.09"
$ws.Range("E3").Value = "  -0.71%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'315.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'104.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.18%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.10%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.601"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.62%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'39.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.79%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0903"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +0.35%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.22%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +3.37%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'15.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.635.31"
$ws.Range("E16").Value = "  -0.59%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.323.78"
$ws.Range("E17").Value = "  +1.39%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.757.13"
$ws.Range("E18").Value = "  +0.72%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'7.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.71%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.95%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "'13.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +22.99%  "

# Row 23 - PancakeSwap
$ws.Range("E23").Value = "  +0.24%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "'262.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.12%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -3.54%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.46%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'10.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28 - Toncoin
$ws.Range("D28").Value = "'2.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.44%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "'7.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +20.91%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "'22.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.20%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "'37.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.15%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'166.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0872"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

# Row 34 - Stellar
$ws.Range("E34").Value = "  -3.87%  "

# Row 35 - WEMIXToken
$ws.Range("D35").Value = "'2.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.73%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -2.76%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  -1.30%  "

# Row 38 - NEARProtocol -> VeChain
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0349"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.92%  "

# Row 39 - VeChain -> NEARProtocol
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  -3.60%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Value = "'1.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.62%  "

# Row 42 - Algorand
$ws.Range("D42").Value = "'0.231"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.44%  "

# Row 43 - MultiversX
$ws.Range("D43").Value = "'69.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

# Row 44 - FirstDigitalUSD
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "

# Row 45 - BitcoinSV
$ws.Range("D45").Value = "'92.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46 - Celestia
$ws.Range("D46").Value = "'12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.00%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'113.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.46%  "

# Row 48 - Maker
$ws.Range("D48").Value = "1.727.88"
$ws.Range("E48").Value = "  +7.94%  "

# Row 49 - ordi
$ws.Range("D49").Value = "'79.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.41%  "

# Row 50 - FraxShare
$ws.Range("D50").Value = "'8.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "

# Row 51 - THORChain
$ws.Range("D51").Value = "'5.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.79%  "

